# Remove file upload functionality
# Append a new row (row 90) of data to each of the four sheets, mirroring
# the pattern already present for row 89 in each sheet (each sheet gets
# one more daily log entry, dated 2025-08-07 11:06:10).

$wb = $excel.ActiveWorkbook

$newRowData = @{
    "MID_LFT_#1" = @{
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"
        D = "0x01,0x1C"
        E = "0x07"
        F = 400
        G = [double]"5.68631262647113e+23"
        H = 284
        I = 7
    }
    "MID_LFT_#2" = @{
        B = "0x01,0x7c"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
        D = "0x01,0x28"
        E = "0x19"
        F = 380
        G = [double]"5.68432987514711e+23"
        H = 296
        I = 25
    }
    "MID_PLT_#1" = @{
        B = "0x00,0x6e"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
        D = "0x00,0x5F"
        E = "0x15"
        F = 110
        G = [double]"5.68631262647113e+23"
        H = 95
        I = 15
    }
    "MID_PLT_#2" = @{
        B = "0x00,0x82"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
        D = "0x00,0x75"
        E = "0x9"
        F = 130
        G = [double]"5.68631262647113e+23"
        H = 117
        I = 9
    }
}

$newRowDate = Get-Date -Year 2025 -Month 8 -Day 7 -Hour 11 -Minute 6 -Second 10

foreach ($sheetName in $newRowData.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 90

    $data = $newRowData[$sheetName]

    $ws.Cells.Item($row, 1).Value = $newRowDate
    $ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($row, 2).Value = $data.B
    $ws.Cells.Item($row, 3).Value = $data.C
    $ws.Cells.Item($row, 4).Value = $data.D
    $ws.Cells.Item($row, 5).Value = $data.E
    $ws.Cells.Item($row, 6).Value = $data.F
    $ws.Cells.Item($row, 7).Value = $data.G
    $ws.Cells.Item($row, 8).Value = $data.H
    $ws.Cells.Item($row, 9).Value = $data.I
}
